# Add a new "TSIDE" column to the NitroXBots sheet (28th test case new steps).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXBots")
$ws.Activate()

# Insert a new column before K, shifting K:AR to L:AS.
$ws.Columns("K:K").Insert()

# Header for the new column.
$ws.Range("K1").Value = "TSIDE"

# New data point added for the 28th test case (row 6).
$ws.Range("K6").Value = "SELL"

# Move the active selection to reflect the newly added steps.
$ws.Range("L10").Select() | Out-Null
